# "altered the names on the plot"
#
# The chart's category axis pulled its labels from Sheet1!M18:M26, which
# held the raw OpenSim muscle/actuator names (plantarflex_r, glmed_r, ...).
# Replace those with human-readable display names so the bar chart shows
# friendly labels instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M18").Value = "plantarflexors"
$ws.Range("M19").Value = "Glut. Med."
$ws.Range("M20").Value = "Hip Flexors"
$ws.Range("M21").Value = "Hip Adductors"
$ws.Range("M22").Value = "Glut. Min."
$ws.Range("M23").Value = "Glut. Max."
$ws.Range("M24").Value = "Hamstrings"
$ws.Range("M25").Value = "Quads"
$ws.Range("M26").Value = "Dorsiflexors"

# Leave the cursor where the author left it after typing the last label.
$ws.Range("M27").Select() | Out-Null
